$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values are plain decimal numbers (e.g. "529.05").
# Excels normal type inference on Range.Value would silently convert those
# to numeric cells, but the source data stores every Price/Volume cell as text
# (note values like "57.859.28" that are not valid numbers at all). Flip the
# handful of purely-numeric-looking cells to Text format before writing them,
# then clear the formatting again so no stray number-format style sticks to
# the cell (matching the unstyled cells in the original sheet).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '57.859.28'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '3.137.33'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '529.05'
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('D6').Value = '138.87'
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.137.66'
$ws.Range('E8').Value = '  +1.20%  '
$ws.Range('D9').Value = '0.448'
$ws.Range('E9').Value = '  +3.52%  '
$ws.Range('E10').Value = '  -0.72%  '
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('E12').Value = '  +3.27%  '
$ws.Range('D13').Value = '3.677.20'
$ws.Range('E13').Value = '  +1.09%  '
$ws.Range('E14').Value = '  +2.55%  '
$ws.Range('E15').Value = '  -2.71%  '
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').Value = '57.999.99'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').Value = '3.137.19'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D19').Value = '6.03'
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('D20').Value = '12.84'
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('D22').Value = '354.36'
$ws.Range('E22').Value = '  +5.43%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '68.72'
$ws.Range('E24').Value = '  +3.27%  '
$ws.Range('E25').Value = '  -0.53%  '
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').Value = '0.0₃0918'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E29').Value = '  +3.91%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('E31').Value = '  -4.79%  '
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').Value = '21.21'
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('E34').Value = '  -0.62%  '
$ws.Range('E35').Value = '  +7.30%  '
$ws.Range('D36').Value = '158.57'
$ws.Range('E36').Value = '  +1.76%  '
$ws.Range('D37').Value = '6.18'
$ws.Range('E37').Value = '  +1.68%  '
$ws.Range('D38').Value = '26.60'
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D40').Value = '0.0670'
$ws.Range('E40').Value = '  +1.39%  '
$ws.Range('E41').Value = '  +5.92%  '
$ws.Range('D42').Value = '4.16'
$ws.Range('E42').Value = '  +6.83%  '
$ws.Range('E43').Value = '  +2.43%  '
$ws.Range('D44').Value = '3.176.67'
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0272'
$ws.Range('E45').Value = '  +5.41%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '36.59'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('D48').Value = '2.325.83'
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').Value = '20.51'
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('D51').Value = '6.03'
$ws.Range('E51').Value = '  +0.63%  '

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
